$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value (from the Dec 17 2023 crypto data refresh)
$changes = [ordered]@{
    "D2" = "41.948.76"
    "E2" = "  -0.85%  "
    "D3" = "2.216.36"
    "E3" = "  -1.46%  "
    "E4" = "  +0.05%  "
    "D5" = "241.57"
    "E5" = "  -2.41%  "
    "E6" = "  -0.46%  "
    "D7" = "73.48"
    "E7" = "  -1.90%  "
    "D9" = "0.610"
    "E9" = "  -0.94%  "
    "D10" = "43.32"
    "E10" = "  +2.09%  "
    "E11" = "  +0.94%  "
    "D12" = "7.11"
    "E12" = "  -1.30%  "
    "E13" = "  -0.05%  "
    "D14" = "2.545.56"
    "E14" = "  -1.50%  "
    "D15" = "14.26"
    "E15" = "  -1.79%  "
    "D16" = "0.843"
    "E16" = "  -1.65%  "
    "D17" = "2.210.92"
    "E17" = "  -1.82%  "
    "D18" = "41.832.60"
    "E18" = "  -0.77%  "
    "E19" = "  +9.52%  "
    "D20" = "72.83"
    "E20" = "  +1.14%  "
    "D21" = "6.15"
    "E21" = "  -0.06%  "
    "D22" = "10.53"
    "E22" = "  +18.64%  "
    "D23" = "229.61"
    "E23" = "  -0.44%  "
    "E24" = "  -5.01%  "
    "E25" = "  +0.08%  "
    "D26" = "11.51"
    "E26" = "  +1.90%  "
    "D27" = "3.57"
    "E27" = "  -1.38%  "
    "E28" = "  -1.90%  "
    "E29" = "  -0.47%  "
    "D30" = "167.21"
    "E30" = "  -1.18%  "
    "D31" = "20.59"
    "E31" = "  -0.52%  "
    "D32" = "5.56"
    "E32" = "  +6.39%  "
    "D33" = "0.0794"
    "E33" = "  -3.71%  "
    "E34" = "  -0.37%  "
    "D35" = "29.19"
    "E35" = "  -4.83%  "
    "D36" = "0.110"
    "E36" = "  -8.21%  "
    "D37" = "4.27"
    "E37" = "  -6.00%  "
    "E38" = "  -1.31%  "
    "D39" = "12.82"
    "E39" = "  -5.67%  "
    "D40" = "65.63"
    "E40" = "  +6.18%  "
    "E41" = "  -3.39%  "
    "D42" = "5.62"
    "E42" = "  -3.42%  "
    "D43" = "0.200"
    "E43" = "  -1.63%  "
    "D44" = "8.71"
    "E44" = "  +0.30%  "
    "D45" = "104.01"
    "E45" = "  -4.00%  "
    "D46" = "0.101"
    "E46" = "  -1.33%  "
    "E47" = "  +5.65%  "
    "E48" = "  -1.02%  "
    "E49" = "  -0.36%  "
    "E50" = "  +0.06%  "
    "D51" = "2.416.72"
    "E51" = "  -1.65%  "
}

# A handful of the new "Price" values (column D) are plain decimal numbers
# (e.g. "241.57") that Excel would otherwise auto-convert to a numeric type
# when assigned via .Value. The source data stores every Price/Volume cell
# as text, so those specific cells are temporarily switched to a text
# number format before the value is written, then restored to the default
# "Normal" style so no stray formatting is left behind.
$textRiskCells = @("D5", "D7", "D9", "D10", "D12", "D15", "D16", "D20", "D21", "D22", "D23", "D26", "D27", "D30", "D31", "D32", "D33", "D35", "D36", "D37", "D39", "D40", "D42", "D43", "D44", "D45", "D46")

foreach ($cellRef in $changes.Keys) {
    $range = $ws.Range($cellRef)
    $needsTextFormat = $textRiskCells -contains $cellRef
    if ($needsTextFormat) {
        $range.NumberFormat = "@"
    }
    $range.Value = $changes[$cellRef]
    if ($needsTextFormat) {
        $range.Style = "Normal"
    }
}
